# Restore revision: update rule-table row values.
#  - C10: was 18, now 1 (the literal shown in the source revision as 1.0)
#  - B11: was the stray shared string "1", now correctly "R40"
# Once "1" is no longer referenced by any cell, it drops out of the
# shared-strings table automatically (handled by the save pipeline).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C10").Value = 1.0
$ws.Range("B11").Value = "R40"
